$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 64

$ws.Cells.Item($row, 2).Value = 54
$ws.Cells.Item($row, 6).Value = "https://dev.to/rahulmishra05/deadlock-avoidance-banker-s-algorithm-operating-system-m04-p05-4fgk"
$ws.Cells.Item($row, 3).Value = "Deadlock Avoidance (Banker's Algorithm) | Operating System - M04 P05"
$ws.Cells.Item($row, 4).Value = Get-Date -Year 2020 -Month 12 -Day 11
$ws.Cells.Item($row, 4).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item($row, 5).Value = "https://programmingport.hashnode.dev/deadlock-avoidance-bankers-algorithm-or-operating-system-m04-p05"

# Extend the table range to include the new row
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("B10:F64"))

# Update the selection / active cell to mirror the source edit
$ws.Range("E64").Select()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 4
